$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data between row 3 and row 4 for columns D, K, L, M, O, P
$d3 = $ws.Range("D3").Value2
$k3 = $ws.Range("K3").Value2
$l3 = $ws.Range("L3").Value2
$m3 = $ws.Range("M3").Value2
$o3 = $ws.Range("O3").Value2
$p3 = $ws.Range("P3").Value2

$d4 = $ws.Range("D4").Value2
$k4 = $ws.Range("K4").Value2
$l4 = $ws.Range("L4").Value2
$m4 = $ws.Range("M4").Value2
$o4 = $ws.Range("O4").Value2
$p4 = $ws.Range("P4").Value2

$ws.Range("D3").Value2 = $d4
$ws.Range("K3").Value2 = $k4
$ws.Range("L3").Value2 = $l4
$ws.Range("M3").Value2 = $m4
$ws.Range("O3").Value2 = $o4
$ws.Range("P3").Value2 = $p4

$ws.Range("D4").Value2 = $d3
$ws.Range("K4").Value2 = $k3
$ws.Range("L4").Value2 = $l3
$ws.Range("M4").Value2 = $m3
$ws.Range("O4").Value2 = $o3
$ws.Range("P4").Value2 = $p3
